$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_PAR"

$values = @(
    0,
    163.73628366104,
    173.213531615806,
    130.398200101066,
    0,
    223.97635555522,
    222.125412799984,
    95.8041960006402,
    0,
    135.574572872455,
    167.69558196315,
    119.656262653326,
    0,
    217.668115751512,
    238.224471688667,
    165.351858128216
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
